$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 61, shifting existing rows 61:95 down to 62:96.
# This preserves all existing data/formatting for the rows below.
$ws.Rows.Item(61).Insert()

# The newly inserted row 61 is currently blank (except for the D column
# style which Excel carries over). Populate the constant columns by
# copying them from the row directly below (row 62, which holds what
# used to be row 61's original data) and then overwrite the columns
# that differ for this new weekly entry (D, J, K, L, M, P).
$ws.Range("A62:R62").Copy()
$ws.Range("A61").PasteSpecial(-4104)
$excel.CutCopyMode = 0

$ws.Range("D61").Value = 45119
$ws.Range("J61").Value = 1400
$ws.Range("K61").Value = 10000
$ws.Range("L61").Value = 12000
$ws.Range("M61").Value = 11000
$ws.Range("P61").Value = 440
